# Update (Removed Auto Arima)
# Updates the "Amazon Mean Forecast" (D), "Amazon P70 Forecast" (E),
# "Amazon P80 Forecast" (F) and "Amazon P90 Forecast" (G) columns on the
# "Forecast Comparison" sheet for rows 2-17 with recalculated forecast
# values now that the Auto Arima model has been removed from the blend.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$data = @{
    2  = @(10, 12, 15, 21)
    3  = @(9,  11, 15, 21)
    4  = @(9,  10, 14, 21)
    5  = @(9,  11, 16, 23)
    6  = @(9,  10, 14, 22)
    7  = @(9,  11, 15, 22)
    8  = @(9,  10, 15, 22)
    9  = @(8,  10, 14, 22)
    10 = @(8,  10, 14, 21)
    11 = @(8,  9,  13, 20)
    12 = @(8,  9,  13, 20)
    13 = @(8,  10, 14, 22)
    14 = @(8,  9,  13, 21)
    15 = @(7,  9,  13, 20)
    16 = @(8,  9,  13, 20)
    17 = @(7,  8,  12, 19)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 4).Value = $values[0]
    $ws.Cells.Item($row, 5).Value = $values[1]
    $ws.Cells.Item($row, 6).Value = $values[2]
    $ws.Cells.Item($row, 7).Value = $values[3]
}
